# Locale change: Russian/English template labels -> Spanish
# (cambios de idioma a events, route, stops, summary, trips)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("template")

# --- Report header labels (column A) / values (column B) ---
$ws.Range("A2").Value = "Tipo de Reporte:"
$ws.Range("B2").Value = "Paradas"

$ws.Range("A4").Value = "Dispositivo:"
# B4 keeps its jxls placeholder formula text (${device.deviceName}) - untouched

$ws.Range("A5").Value = "Grupo:"
# B5 keeps its jxls placeholder formula text (${device.groupName}) - untouched

$ws.Range("A6").Value = "Periodo:"
# B6 keeps its jxls placeholder formula text (${dateTool.format(...)}) - untouched

# --- Table header row (row 8) ---
$ws.Range("A8").Value = "Inicio"
$ws.Range("B8").Value = "Dirección de Inicio"
$ws.Range("C8").Value = "Kilometraje"
$ws.Range("D8").Value = "Fin"
$ws.Range("E8").Value = "Duración"
$ws.Range("F8").Value = "Horas de Motor"
$ws.Range("G8").Value = "Combistible consumido"

# Row 9 (jxls per-row template formulas) is untouched - content unchanged.

# Restore the active selection to match the edited template's saved view.
$ws.Range("D22").Select() | Out-Null
